# Remove the deprecated duplicate header row (row 2) and make the
# remaining header row (row 1) bold, matching the formatting that the
# duplicate row used to carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the bold formatting that lived on the duplicate header row (row 2)
# up onto the real header row (row 1) before that row disappears.
$ws.Range("A2:D2").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Delete the now-redundant duplicate header row entirely, shifting rows 3-5 up to 2-4.
$ws.Rows(2).Delete()
